$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates (M2:T2)
$ws.Range("M2").Value = 0.8596446666666667
$ws.Range("N2").Value = 2.578934
$ws.Range("O2").Value = 0.05286426382906832
$ws.Range("P2").Value = 0.05286426382906832
$ws.Range("Q2").Value = 0.7161054984500002
$ws.Range("R2").Value = 6.444949486050001
$ws.Range("S2").Value = 0.05286426382906832
$ws.Range("T2").Value = 0.05286426382906832

# Row 3 updates (O3,P3,Q3,S3,T3)
$ws.Range("O3").Value = 0.6417658132713033
$ws.Range("P3").Value = 0.6417658132713032
$ws.Range("Q3").Value = 8.693434738575
$ws.Range("S3").Value = 0.6417658132713033
$ws.Range("T3").Value = 0.6417658132713032

# Row 4 updates (O4,P4,S4,T4)
$ws.Range("O4").Value = 0.3053699228996285
$ws.Range("P4").Value = 0.3053699228996284
$ws.Range("S4").Value = 0.3053699228996285
$ws.Range("T4").Value = 0.3053699228996284
